# Add 'My Services' feature row to the API endpoint list sheet.
#
# A new API endpoint row (service/my-services) is inserted above the
# existing "service/:id" row (row 14), pushing all following rows down
# by one. The new row copies the banding/border formatting used by the
# row that currently sits at the bottom of that table section (the
# thick-bottom-border row), matching the style indices the author used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

# Insert a new blank row at position 14; everything below shifts down.
[void]$ws.Rows(14).Insert()

# The row that used to be 19 (last row of the "service/*" table, with the
# thick bottom border) is now row 20. Copy its formatting into the new
# row 14 so the new entry picks up the same borders/shading.
[void]$ws.Range("A20:E20").Copy()
[void]$ws.Range("A14:E14").PasteSpecial(-4122)

# Fill in the new endpoint's data.
$ws.Cells.Item(14, 1).Value = "service/my-services"
$ws.Cells.Item(14, 2).Value = "GET"
$ws.Cells.Item(14, 3).Value = "A felhasználó szolgáltatásainak lekérése"
# D14 (parameter data) and E14 (needs body) stay empty, same as the
# template row.

# Clear the clipboard marching ants / leftover copy mode.
$excel.CutCopyMode = 0

# Match the author's final selection/cursor position.
[void]$ws.Range("D14").Select()
